$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column H, re-using the same style as the other headers (e.g. G1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# "Save" column values for rows 2-12
$saveValues = @(0, 0, 1, 0, 0, 0, 0, 0, 0, 1, 0)

for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
